# Update the time-slot labels in column C (rows 2,3,6,7).
# Rows 4 and 5 already contain the correct values and are left untouched.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C2").Value = "9:30-9:35"
$ws.Range("C3").Value = "9:35-9:40"
$ws.Range("C6").Value = "19:25-19:30"
$ws.Range("C7").Value = "19:30-19:35"

# Move the active selection from B11 to C11, matching the saved cursor
# position recorded in the sheet view.
$ws.Range("C11").Select()
